$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final Grade column (F): average of Midterm Exam, Midterm Paper, Final Exam, Final Paper, rounded
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Formula = "=ROUND(AVERAGE(B2:E2), 0)"
$ws.Range("F3:F25").ClearFormats()
$ws.Range("F3:F25").Formula = "=ROUND(AVERAGE(B3:E3), 0)"

# Pass/Fail column (G): PASS if final grade >= 60, otherwise " FAIL"
$ws.Range("G2").Formula = '=IF(F2>=60,"PASS"," FAIL")'
$ws.Range("G3:G25").Formula = '=IF(F3>=60,"PASS"," FAIL")'

# Letter Grade column (H): A/B/C/D/F based on final grade thresholds
$ws.Range("H2").Formula = '=IF(F2>=90,"A",IF(F2>=80,"B",IF(F2>70,"C",IF(F2>=60,"D",IF(F2<60,"F")))))'
$ws.Range("H3:H25").Formula = '=IF(F3>=90,"A",IF(F3>=80,"B",IF(F3>70,"C",IF(F3>=60,"D",IF(F3<60,"F")))))'

# Update the sheet's selection to reflect what was last worked on
$ws.Range("H2:H25").Select()
